# Rebuild the stations_metadata table with a new column layout and
# several additional input rows (river gauging stations).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the old table body (keeps header row formatting / style).
$ws.Range("A1:L8").ClearContents()

# New header row (columns were reordered and two columns were added:
# station_kodas, and roughness_n moved next to the coordinates).
$headers = @("river_name","station_name","station_code","station_id","station_kodas","x_coord","y_coord","roughness_n","basin_name","datum_offset_cm","min_level_cm","max_level_cm")
for ($c = 1; $c -le $headers.Length; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# Data rows. station_code (col 3) looks numeric ("101", "102", ... "769")
# but must stay text, so force the Text number format on that column
# before writing the value (mirrors typing '101 into Excel).
$ws.Range("C2:C8").NumberFormat = "@"

$rows = @(
    @("Merkys",               "Puvočiai",            "101", 5101, "5101LT", 575000, 5998000, 0.04,  "Nemunas-Merkys", 0,     50,   850),
    @("Nemunas",               "Druskininkai",        "102", 5102, "5102LT", 568500, 5992000, 0.038, "Nemunas-Main",   0,     20,   1000),
    @("Verknė",                "Verbyliškės",         "103", 5103, "5103LT", 521000, 6032000, 0.042, "Nemunas-Verkne", 0,     30,   700),
    @("Nemunas",               "Nemunaičiai",         "104", 5104, "5104LT", 540500, 6042000, 0.038, "Nemunas-Main",   0,     20,   1100),
    @("Merkys",                "Jašiūnai",            "105", 5105, "5105LT", 583500, 6029000, 0.04,  "Nemunas-Merkys", 0,     30,   800),
    @("Šešupė",                "Kudirkos Naumiestis", "106", 5106, "5106LT", 409000, 6022000, 0.045, "Nemunas-Sesupe", 0,     40,   950),
    @("Nemuno atšaka Atmata",  "Rusnė",               "769", 769,  "60004LT",333694, 6132670, 0.03,  "Nemunas-Delta", -1.56, -100,  300)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    for ($c = 1; $c -le $row.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}
